# Generate Report for Handback
# - Flip the "Ready for handoff" status (Overview + per-locale Status columns)
#   to "Handback transform failed" for the c05459db row.
# - Record the handback/handoff file-name mismatch reason in the
#   "Error Detail" column (L) of row 7 on both locale sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B7").Value = "Handback transform failed"
$overview.Range("C7").Value = "Handback transform failed"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = "Handback transform failed"
$zhcn.Range("L7").Value = "Handback file name: 3avnqujt.u0y is different with handoff file name: c05459db-eafd-46b4-84b5-31928066ca5a.319114fd663443f80dd2c6af3176ae741572cb2f.zh-cn."

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = "Handback transform failed"
$dede.Range("L7").Value = "Handback file name: 3avnqujt.u0y is different with handoff file name: c05459db-eafd-46b4-84b5-31928066ca5a.319114fd663443f80dd2c6af3176ae741572cb2f.de-de."
